$wb = $excel.ActiveWorkbook

# --- Expenditures sheet: fix header text (add missing space after '%') ---
$wsExp = $wb.Worksheets.Item("Expenditures")
$wsExp.Range("D1").Value = "% Other  Material Expenditures (collection maintance/binding)"

# --- SpacesStaff sheet: remove the "Seats" and "Study rooms" columns ---
$wsSpaces = $wb.Worksheets.Item("SpacesStaff")
$wsSpaces.Columns("D:E").Delete()

# --- View/selection updates ---
# Expenditures: drop the frozen/scrolled topLeftCell, move the selection.
$wsExp.Activate()
$wsExp.Range("D13").Select()

# Visitors: keep its own selection as-is (no longer the active tab).
$wsVisitors = $wb.Worksheets.Item("Visitors")
$wsVisitors.Range("I10").Select()

# SpacesStaff becomes the active tab, with a new selection.
$wsSpaces.Activate()
$wsSpaces.Range("O13").Select()
